$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Аркуш1" to "ThisYear"
$ws.Name = "ThisYear"

# Move the active selection from P38 to O24
[void]$ws.Range("O24").Select()

# Highlight R9 (the ABS() result cell) with the green fill used elsewhere
# in the column (style index 5: fillId=2 -> FF92D050, centered)
$ws.Range("R9").Interior.Color = 5296274
